$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.600.48'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').Value = '2.289.78'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '96.20'
$ws.Range('E5').Value = '  +3.22%  '
$ws.Range('D6').Value = '266.73'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').Value = '0.623'
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('D9').Value = '0.609'
$ws.Range('E9').Value = '  -1.50%  '
$ws.Range('D10').Value = '45.60'
$ws.Range('E10').Value = '  -0.73%  '
$ws.Range('D11').Value = '0.0934'
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').Value = '7.81'
$ws.Range('E12').Value = '  -2.13%  '
$ws.Range('E13').Value = '  +0.46%  '
$ws.Range('D14').Value = '2.632.81'
$ws.Range('E14').Value = '  +0.14%  '
$ws.Range('D15').Value = '15.13'
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('D16').Value = '0.846'
$ws.Range('E16').Value = '  +1.77%  '
$ws.Range('D17').Value = '2.290.47'
$ws.Range('E17').Value = '  +0.06%  '
$ws.Range('D18').Value = '43.594.22'
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('E19').Value = '  +2.70%  '
$ws.Range('E20').Value = '  +0.54%  '
$ws.Range('D21').Value = '71.94'
$ws.Range('E21').Value = '  +1.88%  '
$ws.Range('E22').Value = '  +7.60%  '
$ws.Range('D23').Value = '232.56'
$ws.Range('E23').Value = '  -0.66%  '
$ws.Range('D24').Value = '9.15'
$ws.Range('E24').Value = '  -8.62%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  +1.82%  '
$ws.Range('D27').Value = '11.14'
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('E28').Value = '  +3.10%  '
$ws.Range('D29').Value = '40.09'
$ws.Range('E29').Value = '  +1.53%  '
$ws.Range('D30').Value = '2.22'
$ws.Range('E30').Value = '  -1.25%  '
$ws.Range('D31').Value = '175.57'
$ws.Range('E31').Value = '  +1.65%  '
$ws.Range('D32').Value = '21.83'
$ws.Range('E32').Value = '  -1.37%  '
$ws.Range('D33').Value = '0.0890'
$ws.Range('E33').Value = '  -2.36%  '
$ws.Range('E34').Value = '  -3.49%  '
$ws.Range('D35').Value = '0.126'
$ws.Range('E35').Value = '  +1.09%  '
$ws.Range('E36').Value = '  -2.13%  '
$ws.Range('D37').Value = '0.0355'
$ws.Range('E37').Value = '  +2.88%  '
$ws.Range('D38').Value = '4.30'
$ws.Range('E38').Value = '  -2.51%  '
$ws.Range('D39').Value = '3.40'
$ws.Range('E39').Value = '  +0.56%  '
$ws.Range('E40').Value = '  -3.14%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').Value = '12.26'
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D44').Value = '64.39'
$ws.Range('E44').Value = '  +6.07%  '
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('E46').Value = '  -3.62%  '
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('D48').Value = '97.56'
$ws.Range('E48').Value = '  -1.52%  '
$ws.Range('E49').Value = '  +1.55%  '
$ws.Range('D50').Value = '2.512.31'
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('D51').Value = '0.429'
$ws.Range('E51').Value = '  +2.25%  '
